$p = $ppt.ActivePresentation
$s = $p.Slides.Add(9, 2)
$s.Shapes.Item(1).Delete()
$s.Shapes.Item(1).Delete()
$r = $s.Shapes.AddShape(1, 100, 100, 200, 50)
$r.Fill.ForeColor.SchemeColor = 9
